$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for electrode "G9" that was missing, right after the
# G7 row (original row 31), shifting everything below down by one.
$ws.Rows(32).Insert()
$ws.Range("A32").Value2 = "G9_bipolar_20V_1kHz.txt"
$ws.Range("B32").Value2 = 59.197821
$ws.Range("C32").Value2 = "G9"

# Replace the "File Name" column (A) with the electrode location code
# that already lives in column C, for every data row.
for ($r = 2; $r -le 69; $r++) {
    $loc = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 1).Value2 = $loc
}

# Update header labels.
$ws.Range("A1").Value2 = "Loc"
$ws.Range("B1").Value2 = "P_max"

# Drop the now-redundant "Electrode Locations" column.
$ws.Columns("C").Delete()
